$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.141.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.637.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.516'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.33%  '

$ws.Range("E7").Value = '  +0.52%  '

$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.867.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.644.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.544'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.157.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0737'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.73%  '

$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("E21").Value = '  +1.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("E26").Value = '  +0.44%  '

$ws.Range("E27").Value = '  -0.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.118'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.38%  '

$ws.Range("E30").Value = '  +0.26%  '

$ws.Range("E31").Value = '  -0.68%  '

$ws.Range("E32").Value = '  +0.76%  '

$ws.Range("E33").Value = '  -0.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.304.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.35%  '

$ws.Range("E35").Value = '  -0.73%  '

$ws.Range("E36").Value = '  +0.99%  '

$ws.Range("E37").Value = '  -0.69%  '

$ws.Range("E38").Value = '  +0.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.853'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.76%  '

$ws.Range("E40").Value = '  +0.36%  '

$ws.Range("E41").Value = '  +3.71%  '

$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.777.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.32%  '

$ws.Range("E46").Value = '  -1.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.54%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.66'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.753'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +13.27%  '
